# Applies updated simulation results for the "380 kV" case (pl_mw.xlsx, Sheet1)
# to columns B,C,D,E,G,H,I,J,K,L,O across rows 2-25, matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => { ColumnIndex => NewValue }
$updates = @{
  2 = @{ 2=0.6461749584804295; 3=0.1094454731220154; 4=0.5322728414686821; 5=0.1790822283264184; 7=1.101306767836768; 8=1.11287034217942; 9=1.146592119600857; 10=0.07093276835435436; 11=0.3744685892183668; 12=0.3979277363066842; 15=4.499097197137942 }
  3 = @{ 2=0.6051873062638151; 3=0.1080092387727802; 4=0.5298351624341961; 5=0.1793319690276789; 7=1.108875982735164; 8=1.120536043491867; 9=1.15628504769483; 10=0.07130996620941765; 11=0.3370222873391242; 12=0.3931358288301539; 15=4.530913041718463 }
  4 = @{ 2=0.5801556332272639; 3=0.1071246504466643; 4=0.5285616155283321; 5=0.179549578067757; 7=1.114082700000225; 8=1.125642704229392; 9=1.162718785918667; 10=0.07155622504394987; 11=0.3140376101979427; 12=0.3903483354083193; 15=4.552460932178505 }
  5 = @{ 2=0.5699897606569948; 3=0.1067635149382085; 4=0.528098913798587; 5=0.1796544541707377; 7=1.116345140688978; 8=1.127824378251645; 9=1.165461925462669; 10=0.07166027066466274; 11=0.3046737066057119; 12=0.3892514826619049; 15=4.561748348516417 }
  6 = @{ 2=0.5683038472794522; 3=0.1067035097456994; 4=0.5280254861700797; 5=0.1796728482476961; 7=1.11672931425737; 8=1.128192727002734; 9=1.165924752492238; 10=0.07167777063293634; 11=0.3031190101247745; 12=0.3890717157054837; 15=4.563321116621339 }
  7 = @{ 2=0.5800183909837244; 3=0.1071197826735215; 4=0.5285551472945684; 5=0.1795509268250814; 7=1.114112642471056; 8=1.125671719296498; 9=1.162755289461067; 10=0.0715576132764939; 11=0.313911314173339; 12=0.3903333844507131; 15=4.552584134381092 }
  8 = @{ 2=0.6320148888684969; 3=0.1089508440385103; 4=0.5313860955303085; 5=0.1791550228523917; 7=1.103800602090161; 8=1.115430546814466; 9=1.149834197667001; 10=0.07105978835723548; 11=0.3615559260043426; 12=0.3962434595377715; 15=4.509649735324814 }
  9 = @{ 2=0.7350199356715166; 3=0.1125187490943063; 4=0.5387025859608343; 5=0.1788870300156162; 7=1.0880139078632; 8=1.098515636720961; 9=1.128318995834992; 10=0.07019952569231602; 11=0.4550223619599763; 12=0.4090549998408903; 15=4.441412735229534 }
  10 = @{ 2=0.8112988037186142; 3=0.1151249632013673; 4=0.5451471737181777; 5=0.1789981321664662; 7=1.079117519343711; 8=1.08801293200959; 9=1.114837850141342; 10=0.0696377251800353; 11=0.5236882681654151; 12=0.4192058211320102; 15=4.400989906185174 }
  11 = @{ 2=0.8461240527647931; 3=0.1163070711777294; 4=0.5483097663072414; 5=0.1791151585722872; 7=1.075656744570111; 8=1.083651522768989; 9=1.109209065870793; 10=0.06939729975464459; 11=0.5549203422300764; 12=0.4239826432295928; 15=4.384705677532395 }
  12 = @{ 2=0.8593288137810191; 3=0.1167541798262377; 4=0.5495404244782947; 5=0.1791689981858759; 7=1.074430508118326; 8=1.082059730203625; 9=1.107149977206344; 10=0.06930842664851866; 11=0.5667459158395047; 12=0.4258142476893028; 15=4.378841578446156 }
  13 = @{ 2=0.8564841812683142; 3=0.1166579110385157; 4=0.5492739126135859; 5=0.1791569798454269; 7=1.074690851796433; 8=1.082399893952925; 9=1.107590219049179; 10=0.06932747061703015; 11=0.564199137532313; 12=0.4254187707458073; 15=4.380091071859283 }
  14 = @{ 2=0.8472100767491213; 3=0.1163438658734961; 4=0.5484103517396051; 5=0.1791193973171872; 7=1.07555417228059; 8=1.083519367491547; 9=1.109038212459588; 10=0.06938994464120718; 11=0.5558932700666901; 12=0.4241328759021172; 15=4.384217175551527 }
  15 = @{ 2=0.8415316335369027; 3=0.1161514342123979; 4=0.5478856957125515; 5=0.1790976161963016; 7=1.076093956757148; 8=1.084212859663168; 9=1.109934578816397; 10=0.06942849429603459; 11=0.5508054917115714; 12=0.4233481829245136; 15=4.386783906300138 }
  16 = @{ 2=0.809025345167754; 3=0.1150476371684803; 4=0.5449451244567598; 5=0.1789918189714932; 7=1.079355484316281; 8=1.088306330292014; 9=1.115215839531473; 10=0.06965374166400906; 11=0.521647030123745; 12=0.4188968335010372; 15=4.40209644484159 }
  17 = @{ 2=0.7891153437629725; 3=0.1143695822663986; 4=0.5432002141406258; 5=0.178943918503343; 7=1.08150646372053; 8=1.090924109398742; 9=1.118584740165936; 10=0.06979579689700177; 11=0.5037576324624808; 12=0.4162067250325379; 15=4.412029005846108 }
  18 = @{ 2=0.7776755122128804; 3=0.1139792580540018; 4=0.542218332041557; 5=0.1789226275162505; 7=1.082798830584139; 8=1.092468978321378; 9=1.120569870693103; 10=0.06987892877979007; 11=0.4934677527029407; 12=0.4146744365968118; 15=4.417940033455039 }
  19 = @{ 2=0.7738042477677709; 3=0.1138470461094201; 4=0.5418896228638772; 5=0.1789164950370363; 7=1.083245881576715; 8=1.092998778295367; 9=1.121250148933651; 10=0.06990732082829698; 11=0.4899837339077635; 12=0.414158209676188; 15=4.419975432994931 }
  20 = @{ 2=0.7912335747109296; 3=0.1144417962260391; 4=0.543383713489149; 5=0.1789483699567604; 7=1.081271777442538; 8=1.090641386799632; 9=1.118221206643597; 10=0.06978052739077878; 11=0.5056620329227997; 12=0.4164915413521584; 15=4.410951168651877 }
  21 = @{ 2=0.8499336462886617; 3=0.1164361231789499; 4=0.5486631046697852; 5=0.1791301780201202; 7=1.075298306820088; 8=1.08318892930869; 9=1.108610936865581; 10=0.06937153564209186; 11=0.558332946507619; 12=0.4245099589110595; 15=4.382997033942388 }
  22 = @{ 2=0.8883973610872431; 3=0.1177364278361352; 4=0.5523060792057777; 5=0.1793044979274292; 7=1.071885561316734; 8=1.078666711458055; 9=1.102752109260187; 10=0.06911688595197862; 11=0.5927483267622051; 12=0.4298828454748929; 15=4.366489822362837 }
  23 = @{ 2=0.8678597046847187; 3=0.1170427248230226; 4=0.5503441838607728; 5=0.1792063940342139; 7=1.073662062744546; 8=1.081048455247199; 9=1.105840473741445; 10=0.06925164190332467; 11=0.5743811643255299; 12=0.4270031748038292; 15=4.375138839838314 }
  24 = @{ 2=0.7902759015847209; 3=0.1144091498538771; 4=0.5433006871067789; 5=0.1789463379914906; 7=1.08137770548079; 8=1.090769081477603; 9=1.118385409679529; 10=0.0697874261788094; 11=0.5048010690796616; 12=0.4163627313685367; 15=4.411437833768446 }
  25 = @{ 2=0.7070464148094402; 3=0.1115561047261338; 4=0.5365350312876558; 5=0.1789053031168351; 7=1.091809985801206; 8=1.102753098205312; 9=1.133730618168585; 10=0.0704198822555977; 11=0.4297360613533385; 12=0.4054590210379132; 15=4.458166086354794 }
}

foreach ($rowKey in $updates.Keys) {
  $rowNum = [int]$rowKey
  $cols = $updates[$rowKey]
  foreach ($colKey in $cols.Keys) {
    $colNum = [int]$colKey
    $ws.Cells.Item($rowNum, $colNum).Value = $cols[$colKey]
  }
}

Write-Host "Updated $($updates.Count) rows across pl_mw.xlsx Sheet1"